$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5276
$wsExpo.Range("F3").Value = 572
$wsExpo.Range("F4").Value = 10712
$wsExpo.Range("F5").Value = 269
$wsExpo.Range("F6").Value = 573
$wsExpo.Range("F7").Value = 149
$wsExpo.Range("F8").Value = 177
$wsExpo.Range("F9").Value = 887

# Sheet "演出" (Performances) - update "想去人数" (F column) values
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 20

# Sheet "全部类型" (All Types) - update "想去人数" (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5276
$wsAll.Range("F5").Value = 572
$wsAll.Range("F6").Value = 20
$wsAll.Range("F7").Value = 10712
$wsAll.Range("F8").Value = 269
$wsAll.Range("F9").Value = 573
$wsAll.Range("F10").Value = 149
$wsAll.Range("F13").Value = 177
$wsAll.Range("F14").Value = 887
